$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 1.02
    "C2" = 1.036038093016071
    "D2" = 1.045322645424396
    "E2" = 1.035001808334363
    "F2" = 1.052172954246519
    "I2" = 1.032677687955853
    "J2" = 1.041148591938893
    "K2" = 1.048091209690816
    "L2" = 1.0377996643153
    "M2" = 1.05492242056042
    "N2" = 1.042627142360859
    "B3" = 1.02
    "C3" = 1.037116611126068
    "D3" = 1.046372156142156
    "E3" = 1.035922098123694
    "F3" = 1.053371257679468
    "I3" = 1.032813050072243
    "J3" = 1.041870257755953
    "K3" = 1.04895151746996
    "L3" = 1.038528963353022
    "M3" = 1.055932527072181
    "N3" = 1.043349833026156
    "B4" = 1.02
    "C4" = 1.03781451815435
    "D4" = 1.04705163237875
    "E4" = 1.036517969542671
    "F4" = 1.054147319026381
    "I4" = 1.032899070497003
    "J4" = 1.042336700454698
    "K4" = 1.049507966197953
    "L4" = 1.039000625054044
    "M4" = 1.056586240975379
    "N4" = 1.043816938127034
    "B5" = 1.02
    "C5" = 1.038107926744809
    "D5" = 1.047337373397861
    "E5" = 1.036768565136205
    "F5" = 1.054473738030309
    "I5" = 1.032934857880261
    "J5" = 1.042532667611929
    "K5" = 1.049741842678528
    "L5" = 1.039198853173767
    "M5" = 1.056861088170872
    "N5" = 1.044013183580107
    "B6" = 1.02
    "C6" = 1.038157191847206
    "D6" = 1.047385355803666
    "E6" = 1.03681064658133
    "F6" = 1.054528554773184
    "I6" = 1.032940844706807
    "J6" = 1.042565564007282
    "K6" = 1.049781108369125
    "L6" = 1.039232133114732
    "M6" = 1.056907237745784
    "N6" = 1.044046126692114
    "B7" = 1.02
    "C7" = 1.037818438660477
    "D7" = 1.0470554501137
    "E7" = 1.036521317655086
    "F7" = 1.054151680012984
    "I7" = 1.032899550165937
    "J7" = 1.042339319470085
    "K7" = 1.049511091481337
    "L7" = 1.039003274017787
    "M7" = 1.056589913395221
    "N7" = 1.043819560861723
    "B8" = 1.02
    "C8" = 1.036402576530892
    "D8" = 1.045677255332853
    "E8" = 1.035312745224713
    "F8" = 1.052577786102696
    "I8" = 1.032723758719595
    "J8" = 1.041392590603611
    "K8" = 1.048382001929579
    "L8" = 1.03804618467442
    "M8" = 1.055263768626857
    "N8" = 1.042871487531668
    "B9" = 1.02
    "C9" = 1.033907885855682
    "D9" = 1.043251553211309
    "E9" = 1.033186027369787
    "F9" = 1.049809570297389
    "I9" = 1.032401993828691
    "J9" = 1.039720327574775
    "K9" = 1.046390655542878
    "L9" = 1.036357812095791
    "M9" = 1.052927747913911
    "N9" = 1.041196849697519
    "B10" = 1.02
    "C10" = 1.03224489053048
    "D10" = 1.041636332648305
    "E10" = 1.031770207337506
    "F10" = 1.047967558538452
    "I10" = 1.032179430436848
    "J10" = 1.038602792960582
    "K10" = 1.045061913406661
    "L10" = 1.035230986847905
    "M10" = 1.05137093638074
    "N10" = 1.040077728055992
    "B11" = 1.02
    "C11" = 1.031524819067613
    "D11" = 1.040937376418306
    "E11" = 1.031157616290795
    "F11" = 1.047170763286696
    "I11" = 1.032081150355034
    "J11" = 1.038118248810326
    "K11" = 1.044486272288781
    "L11" = 1.034742764689572
    "M11" = 1.050696943036918
    "N11" = 1.039592495797464
    "B12" = 1.02
    "C12" = 1.03125735419599
    "D12" = 1.040677819523699
    "E12" = 1.030930143152121
    "F12" = 1.046874918836898
    "I12" = 1.032044358150957
    "J12" = 1.037938170583317
    "K12" = 1.044272410173761
    "L12" = 1.034561372213942
    "M12" = 1.05044660893969
    "N12" = 1.039412161838717
    "B13" = 1.02
    "C13" = 1.031314726252529
    "D13" = 1.040733492350671
    "E13" = 1.030978933723121
    "F13" = 1.046938373019674
    "I13" = 1.032052263170747
    "J13" = 1.037976802364282
    "K13" = 1.04431828629023
    "L13" = 1.034600283554108
    "M13" = 1.050500305689218
    "N13" = 1.039450848481242
    "B14" = 1.02
    "C14" = 1.031502710306889
    "D14" = 1.040915920008012
    "E14" = 1.031138811852517
    "F14" = 1.047146306231744
    "I14" = 1.032078114945776
    "J14" = 1.038103365466682
    "K14" = 1.044468595265915
    "L14" = 1.034727771649789
    "M14" = 1.050676249997254
    "N14" = 1.039577591317765
    "B15" = 1.02
    "C15" = 1.031618533694104
    "D15" = 1.041028328521523
    "E15" = 1.031237327382711
    "F15" = 1.047274436733246
    "I15" = 1.032094005104006
    "J15" = 1.038181332318921
    "K15" = 1.044561199831652
    "L15" = 1.034806315305087
    "M15" = 1.050784657337293
    "N15" = 1.039655668891877
    "B16" = 1.02
    "C16" = 1.03229267958562
    "D16" = 1.041682729422015
    "E16" = 1.031810872863021
    "F16" = 1.048020456186732
    "I16" = 1.032185912761606
    "J16" = 1.038634936972127
    "K16" = 1.045100110743278
    "L16" = 1.035263382232761
    "M16" = 1.051415669512797
    "N16" = 1.04010991771572
    "B17" = 1.02
    "C17" = 1.032715557379952
    "D17" = 1.042093336744608
    "E17" = 1.032170768246286
    "F17" = 1.048488630374356
    "I17" = 1.032243053109147
    "J17" = 1.038919298581573
    "K17" = 1.045438078301927
    "L17" = 1.035550008052424
    "M17" = 1.0518115176148
    "N17" = 1.040394683151268
    "B18" = 1.02
    "C18" = 1.032962216483614
    "D18" = 1.042332879997035
    "E18" = 1.032380734319791
    "F18" = 1.048761786457701
    "I18" = 1.032276197938695
    "J18" = 1.039085099743586
    "K18" = 1.045635181386227
    "L18" = 1.035717163160193
    "M18" = 1.052042420368235
    "N18" = 1.040560719769948
    "B19" = 1.02
    "C19" = 1.033046321164732
    "D19" = 1.042414565348882
    "E19" = 1.032452334986478
    "F19" = 1.048854938871899
    "I19" = 1.032287468226314
    "J19" = 1.039141623117597
    "K19" = 1.045702383757993
    "L19" = 1.035774153834563
    "M19" = 1.052121154162028
    "N19" = 1.040617323413634
    "B20" = 1.02
    "C20" = 1.032670186430203
    "D20" = 1.04204927802437
    "E20" = 1.032132150203218
    "F20" = 1.048438391651037
    "I20" = 1.032236941538915
    "J20" = 1.038888795701935
    "K20" = 1.045401820445607
    "L20" = 1.035519258802372
    "M20" = 1.051769045698798
    "N20" = 1.040364136954042
    "B21" = 1.02
    "C21" = 1.031447353646464
    "D21" = 1.040862197774367
    "E21" = 1.031091729782787
    "F21" = 1.047085071754037
    "I21" = 1.032070510152892
    "J21" = 1.038066098465242
    "K21" = 1.04442433421829
    "L21" = 1.03469023082655
    "M21" = 1.050624438316939
    "N21" = 1.039540271392908
    "B22" = 1.02
    "C22" = 1.030678519787149
    "D22" = 1.040116218324248
    "E22" = 1.03043798382259
    "F22" = 1.046234884340409
    "I22" = 1.031964209995899
    "J22" = 1.037548275027623
    "K22" = 1.043809498759626
    "L22" = 1.034168727491885
    "M22" = 1.049904876600608
    "N22" = 1.039021712586613
    "B23" = 1.02
    "C23" = 1.031086092497378
    "D23" = 1.04051163974217
    "E23" = 1.030784508256999
    "F23" = 1.046685518765877
    "I23" = 1.032020718843043
    "J23" = 1.037822836261935
    "K23" = 1.044135458553036
    "L23" = 1.034445210989523
    "M23" = 1.050286320701864
    "N23" = 1.03929666372937
    "B24" = 1.02
    "C24" = 1.032690687604062
    "D24" = 1.042069186130657
    "E24" = 1.032149599894852
    "F24" = 1.048461092126322
    "I24" = 1.032239703663453
    "J24" = 1.03890257883215
    "K24" = 1.0454182038955
    "L24" = 1.035533153153999
    "M24" = 1.051788236894456
    "N24" = 1.040377939657883
    "B25" = 1.02
    "C25" = 1.034552797833865
    "D25" = 1.043878317659246
    "E25" = 1.033735484876348
    "F25" = 1.050524608049735
    "I25" = 1.032486598210636
    "J25" = 1.040153122284645
    "K25" = 1.046905673647436
    "L25" = 1.03679451651061
    "M25" = 1.053531569912094
    "N25" = 1.04163025902552
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

Write-Host "Updated $($values.Count) cells"